$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Obrigatorio" column (E) for rows 2 through 8 from "N" to "S"
for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 5).Value = "S"
}
